$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lorem = "Lorem ipsum dolor sit amet, consectetur adipiscing elit. Aenean mollis ultricies interdum. Nullam pharetra vitae lectus eget volutpat. Integer in sodales ligula. Vestibulum pellentesque arcu in est aliquam rhoncus. Curabitur et dui quis arcu scelerisque congue. Pellentesque libero ligula, sagittis a tempus quis, finibus eget erat. Nunc sed tempor nunc. Mauris tempor odio id lorem commodo dapibus. Nulla viverra mi in magna imperdiet volutpat."

# Fix the tiny floating point drift on the existing row 4 timestamp.
$ws.Range("A4").Value = 45778.75077079861

# Grab the date/time number format used by the existing timestamp column
# so newly appended rows stay consistent with the rest of column A.
$dateFormat = $ws.Range("A4").NumberFormat

# New submission rows (timestamp, submissionid, email, phone, firstName, lastName, feedbackText)
$rows = @(
    @{ R=5;  A=45789.37653412037;   B="e7a95397-2efe-11f0-8dd5-fa163ee583d0"; E="john@example.com";             F="";       G="John"; H="Smith";  I=$lorem },
    @{ R=6;  A=45789.3798040162;    B="aeb3eab1-2eff-11f0-8dd5-fa163ee583d0"; E="john@example.com";             F="";       G="John"; H="Smith";  I=$lorem },
    @{ R=7;  A=45789.3834979051;    B="6a36d908-2f00-11f0-8dd5-fa163ee583d0"; E="john@example.com";             F="+41312"; G="John"; H="Smith";  I=$lorem },
    @{ R=8;  A=45789.4085253125;    B="74fdb053-2f05-11f0-8dd5-fa163ee583d0"; E="kevin.maier@students.fhnw.ch"; F="312312"; G="John"; H="Smith";  I=$lorem },
    @{ R=9;  A=45789.41376583333;   B="8348a2fa-2f06-11f0-8dd5-fa163ee583d0"; E="john@example.com";             F="+41312"; G="John"; H="Smith";  I=$lorem },
    @{ R=10; A=45789.41448775463;   B="a8f0fcc1-2f06-11f0-8dd5-fa163ee583d0"; E="kevin.maier@students.fhnw.ch"; F="312312"; G="John"; H="Smith";  I=$lorem },
    @{ R=11; A=45789.41756166667;   B="44d6f554-2f07-11f0-8dd5-fa163ee583d0"; E="john@example.com";             F="+41312"; G="John"; H="Smith";  I=$lorem },
    @{ R=12; A=45789.41876859953;   B="835e7f6f-2f07-11f0-8dd5-fa163ee583d0"; E="kevin.maier@students.fhnw.ch"; F="312312"; G="John"; H="Smith";  I=$lorem },
    @{ R=13; A=45789.45248675926;   B="4d646cfe-2f0e-11f0-8dd5-fa163ee583d0"; E="kevin.maier@students.fhnw.ch"; F="312312"; G="John"; H="Smith";  I=$lorem },
    @{ R=14; A=45789.47321222222;   B="76ca6e57-2f12-11f0-8dd5-fa163ee583d0"; E="kevin.maier@students.fhnw.ch"; F="312312"; G="John"; H="Smith";  I=$lorem },
    @{ R=15; A=45789.51575185185;   B="06adef80-2f1b-11f0-8dd5-fa163ee583d0"; E="kevin.maier@students.fhnw.ch"; F="312312"; G="John"; H="Smith";  I=$lorem },
    @{ R=16; A=45789.51829799668;   B="8ae1c6ae-2f1b-11f0-8dd5-fa163ee583d0"; E="kevin.maier@students.fhnw.ch"; F="312312"; G="Olaf"; H="Schulz"; I=$lorem }
)

foreach ($row in $rows) {
    $r = $row.R

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 1).NumberFormat = $dateFormat

    $ws.Cells.Item($r, 2).Value = $row.B

    $ws.Cells.Item($r, 5).Value = $row.E

    if ($row.F -ne "") {
        # Phone values such as "+41312"/"312312" look numeric to Excel's
        # auto-detection, which would strip the leading "+" or convert the
        # value to a number. Briefly force a text format while assigning,
        # then restore the "Normal" style so the cell keeps plain text
        # without leaving a custom number format applied to it.
        $ws.Cells.Item($r, 6).NumberFormat = "@"
        $ws.Cells.Item($r, 6).Value = $row.F
        $ws.Cells.Item($r, 6).Style = "Normal"
    }

    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
}
